$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'youth knee sleeves with pads'
$ws.Cells.Item(2, 1).Value = 'youth gloves and knee pads'
$ws.Cells.Item(3, 1).Value = 'mens tights 3 4'
$ws.Cells.Item(4, 1).Value = 'black leggings capri'
$ws.Cells.Item(5, 1).Value = 'knee pad tights'
$ws.Cells.Item(6, 1).Value = 'knee leggings'
$ws.Cells.Item(7, 1).Value = 'mens compression pants green'
$ws.Cells.Item(8, 1).Value = 'men tights leggings'
$ws.Cells.Item(9, 1).Value = 'knee pad workout'
$ws.Cells.Item(10, 1).Value = 'hex knee pads basketball'
$ws.Cells.Item(11, 1).Value = 'white knee pads basketball'
$ws.Cells.Item(12, 1).Value = 'knee pad basketball leggings'
$ws.Cells.Item(13, 1).Value = 'capri compression tights'
$ws.Cells.Item(14, 1).Value = 'compression workout capri'
$ws.Cells.Item(15, 1).Value = 'knee pad leggings youth'
$ws.Cells.Item(16, 1).Value = 'all white knee pads for basketball'
$ws.Cells.Item(17, 1).Value = 'basketball compression pants with knee pads youth'
$ws.Cells.Item(18, 1).Value = 'triple 8 knee pads youth'
$ws.Cells.Item(19, 1).Value = 'mens leggings compression under armour'
$ws.Cells.Item(20, 1).Value = 'dirt bike knee pads youth'
$ws.Cells.Item(21, 1).Value = 'nike leggings men basketball'
$ws.Cells.Item(22, 1).Value = 'gym shark mens leggings'
$ws.Cells.Item(23, 1).Value = 'multicam pants with knee pads'
$ws.Cells.Item(24, 1).Value = 'tactical pants with knee pads for men'
$ws.Cells.Item(25, 1).Value = 'excersize gear for men'
$ws.Cells.Item(26, 1).Value = 'men capri leggings'
$ws.Cells.Item(27, 1).Value = 'mens leggings black'
$ws.Cells.Item(28, 1).Value = 'under armour tights'
$ws.Cells.Item(29, 1).Value = 'guys tights'
$ws.Cells.Item(30, 1).Value = 'leggings for basketball'
$ws.Cells.Item(31, 1).Value = 'dri fit compression pants men'
$ws.Cells.Item(32, 1).Value = 'pink compression pants men'
$ws.Cells.Item(33, 1).Value = 'men workout compression pants'
$ws.Cells.Item(34, 1).Value = 'mens football gear'
$ws.Cells.Item(35, 1).Value = 'compression pants with pads basketball'
$ws.Cells.Item(36, 1).Value = 'mens white leggings compression'
$ws.Cells.Item(37, 1).Value = 'cold gear mens'
$ws.Cells.Item(38, 1).Value = 'basketball leggings with kneepads'
$ws.Cells.Item(39, 1).Value = 'mens athletic leggins'
$ws.Cells.Item(40, 1).Value = 'mens compression tights basketball'
$ws.Cells.Item(41, 1).Value = 'under armour cold gear compression pants men'
$ws.Cells.Item(42, 1).Value = 'green mens compression leggings'
$ws.Cells.Item(43, 1).Value = 'mens compression tights leggings'
$ws.Cells.Item(44, 1).Value = 'compression knee pads for basketball'
$ws.Cells.Item(45, 1).Value = 'legging basketball men'
$ws.Cells.Item(46, 1).Value = 'basketball compression pants with padded knees'
$ws.Cells.Item(47, 1).Value = 'knee pad pants'
$ws.Cells.Item(48, 1).Value = 'kneepad honeycomb'
$ws.Cells.Item(49, 1).Value = 'black basketball knee pads'
$ws.Cells.Item(50, 1).Value = 'cheap knee pads for basketball'
$ws.Cells.Item(51, 1).Value = 'compression pants'
$ws.Cells.Item(52, 1).Value = 'knee pads pants'
$ws.Cells.Item(53, 1).Value = 'basketball youth compression pants'
$ws.Cells.Item(54, 1).Value = 'knee pad hex'
$ws.Cells.Item(55, 1).Value = 'men tights and leggings'
$ws.Cells.Item(56, 1).Value = 'knee pad lacrosse'
$ws.Cells.Item(57, 1).Value = 'mens compression tights'
$ws.Cells.Item(58, 1).Value = 'patella guard'
$ws.Cells.Item(59, 1).Value = 'squat pad knee'
$ws.Cells.Item(60, 1).Value = 'best knee pads basketball'
$ws.Cells.Item(61, 1).Value = 'kneepad basketball'
$ws.Cells.Item(62, 1).Value = 'compression mens running pants'
$ws.Cells.Item(63, 1).Value = 'knee protection soccer'
$ws.Cells.Item(64, 1).Value = 'sliding leg guard'
$ws.Cells.Item(65, 1).Value = 'athletic leggings youth'
$ws.Cells.Item(66, 1).Value = 'basketball tights boys'
$ws.Cells.Item(67, 1).Value = 'knee pads xxl'
$ws.Cells.Item(68, 1).Value = 'wrestling clothes for boys'
$ws.Cells.Item(69, 1).Value = 'boys xxl baseball pants'
$ws.Cells.Item(70, 1).Value = 'knee pads girls volleyball'
$ws.Cells.Item(71, 1).Value = 'youth boys leggings sports'
$ws.Cells.Item(72, 1).Value = 'mens compression running pants'
$ws.Cells.Item(73, 1).Value = 'black knee pads wrestling'
$ws.Cells.Item(74, 1).Value = 'black knee pads youth'
$ws.Cells.Item(75, 1).Value = 'hex pad'
$ws.Cells.Item(76, 1).Value = 'men spandex pants'
$ws.Cells.Item(77, 1).Value = 'knee protectors for men'
$ws.Cells.Item(78, 1).Value = 'volleyball gear'
$ws.Cells.Item(79, 1).Value = 'capris for men'
$ws.Cells.Item(80, 1).Value = 'youth black baseball pants'
$ws.Cells.Item(81, 1).Value = 'knee pads professional'
$ws.Cells.Item(82, 1).Value = 'mens mesh pants'
$ws.Cells.Item(83, 1).Value = 'knee pad work pants'
$ws.Cells.Item(84, 1).Value = 'knee pads small'
$ws.Cells.Item(85, 1).Value = 'youth leggings sports'
$ws.Cells.Item(86, 1).Value = 'fitness knee pads'
$ws.Cells.Item(87, 1).Value = 'knee pad sports'
$ws.Cells.Item(88, 1).Value = 'basketball clothes'
$ws.Cells.Item(89, 1).Value = 'girls workout leggings'
$ws.Cells.Item(90, 1).Value = 'baseball pants mens'
$ws.Cells.Item(91, 1).Value = 'kids basketball leggings with knee pads'
$ws.Cells.Item(92, 1).Value = 'youth asics wrestling knee pads'
$ws.Cells.Item(93, 1).Value = 'tesla wintergear for men'
$ws.Cells.Item(94, 1).Value = 'mcdavid basketball knee pads white'
$ws.Cells.Item(95, 1).Value = 'thermo ball mens'
$ws.Cells.Item(96, 1).Value = 'mens nike thermal training pants'
$ws.Cells.Item(97, 1).Value = 'man winter leggings'
$ws.Cells.Item(98, 1).Value = 'underware pants men'
$ws.Cells.Item(99, 1).Value = 'under armor youth basketball compression pants'
$ws.Cells.Item(100, 1).Value = 'black capris'